$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.822.76'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.09%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.639.10'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.38%  '
$ws.Range('E4').Value = '  -0.61%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.92'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.82%  '
$ws.Range('E6').Value = '  -0.30%  '
$ws.Range('E7').Value = '  -0.46%  '
$ws.Range('E8').Value = '  -0.37%  '
$ws.Range('E9').Value = '  -0.60%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.33'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.59%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0845'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.47%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.867.74'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.27%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.628.54'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.71%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.15'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.75%  '
$ws.Range('E15').Value = '  -0.27%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.95'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.32%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.810.21'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.03%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0733'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.75%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '216.25'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.75%  '
$ws.Range('E20').Value = '  -0.44%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.36'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.56'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.23%  '
$ws.Range('E23').Value = '  -2.58%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.17'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.20%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.14'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.45%  '
$ws.Range('E26').Value = '  -0.45%  '
$ws.Range('E27').Value = '  +0.21%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.07'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.49%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.75'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0506'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.84%  '
$ws.Range('E31').Value = '  +1.35%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.37'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.31%  '
$ws.Range('E33').Value = '  -0.65%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.55'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.61%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.263.04'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.92%  '
$ws.Range('E36').Value = '  -0.02%  '
$ws.Range('E37').Value = '  -0.19%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.531'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.79%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.818'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.26%  '
$ws.Range('E40').Value = '  -0.43%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.807'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.39%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.34'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.41%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.779.18'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.77%  '
$ws.Range('E44').Value = '  -4.67%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '92.25'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.66%  '
$ws.Range('E46').Value = '  +0.24%  '
$ws.Range('E47').Value = '  -2.06%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0₆0102'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.22%  '
$ws.Range('E49').Value = '  -1.16%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.56'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0964'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.50%  '
